$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo "ecvellent" -> "excellent" in cell E11
$ws.Range("E11").Value = "excellent"

# Update the selected cell to match the edited cell
$ws.Range("E11").Select()
